# Nanobody_Library_ConditionTest_2 / Growth_Curves workbook cleanup.
#
# The raw-OD600 columns CG, CH, CR, CS and CT (rows 2-42, i.e. the live
# timecourse data under the header row) were being pulled from the wrong
# wells and plotted blank/empty wells that should instead have mirrored the
# no-growth control already recorded in column CF for that timepoint.
# Overwrite each bad reading with the correct value from column CF on the
# same row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nanobody_Slay_test_2.txt")
$ws.Activate()

# Column numbers (1-based): CF=84, CG=85, CH=86, CR=96, CS=97, CT=98
$sourceCol = 84
$targetCols = @(85, 86, 96, 97, 98)

for ($r = 2; $r -le 42; $r++) {
    $correct = $ws.Cells.Item($r, $sourceCol).Value2
    foreach ($c in $targetCols) {
        $ws.Cells.Item($r, $c).Value = $correct
    }
}

# Leave the sheet selection where the author ended up after the fix.
$ws.Range("CT2:CT42").Select() | Out-Null
